$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the header style from H1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row data: row number, I (I0) value, J (IF) value
$rowData = @(
    @(2,6,6),
    @(3,8,8),
    @(4,8,8),
    @(5,8,8),
    @(6,8,8),
    @(7,9,9),
    @(8,8,8),
    @(9,7,7),
    @(10,7,8),
    @(11,9,9),
    @(12,7,8),
    @(13,9,9),
    @(14,8,8),
    @(15,7,8),
    @(16,8,8),
    @(17,8,8),
    @(18,8,8),
    @(19,7,7),
    @(20,8,8),
    @(21,9,9),
    @(22,8,8),
    @(23,8,8),
    @(24,8,8),
    @(25,8,8),
    @(26,8,8),
    @(27,8,8),
    @(28,8,8),
    @(29,8,8),
    @(30,8,9),
    @(31,8,8),
    @(32,8,8),
    @(33,8,8),
    @(34,8,8),
    @(35,9,9),
    @(36,8,8),
    @(37,8,8),
    @(38,7,7),
    @(39,8,8),
    @(40,8,8),
    @(41,7,7),
    @(42,8,8),
    @(43,8,8),
    @(44,8,8),
    @(45,8,8),
    @(46,8,8),
    @(47,8,8),
    @(48,8,8),
    @(49,8,8),
    @(50,7,7),
    @(51,8,8),
    @(52,7,8),
    @(53,7,7),
    @(54,8,8),
    @(55,9,9),
    @(56,9,9),
    @(57,7,7),
    @(58,10,10),
    @(59,8,8),
    @(60,7,8),
    @(61,8,8),
    @(62,9,9),
    @(63,8,8),
    @(64,8,8),
    @(65,8,8),
    @(66,8,8),
    @(67,8,8),
    @(68,9,9),
    @(69,8,8),
    @(70,8,8),
    @(71,9,9),
    @(72,9,9),
    @(73,8,8),
    @(74,5,6),
    @(75,8,8),
    @(76,8,8),
    @(77,8,8),
    @(78,6,6),
    @(79,7,7),
    @(80,7,7),
    @(81,8,8),
    @(82,6,6),
    @(83,4,4),
    @(84,5,6),
    @(85,6,6),
    @(86,6,6),
    @(87,5,5)
)

foreach ($item in $rowData) {
    $r = $item[0]
    $ws.Cells.Item($r, 9).Value = $item[1]
    $ws.Cells.Item($r, 10).Value = $item[2]
}
